$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "VALOR MORA" total (E11) ---
$ws.Range("E11").Value = 416267

# --- 2. Update "Cant. Periodos" (F13) ---
$ws.Range("F13").Value = 8

# --- 3. Insert a new data row (row 23) below the existing last data row (row 22),
#        copying format+values from row 22 into the new row 23 first ---
$ws.Rows.Item(23).Insert(-4121)   # xlShiftDown
$ws.Range("B22:J22").Copy($ws.Range("B23:J23"))
$excel.CutCopyMode = $false

# Reset row 22 back to the "regular" (non-last) row style, copied from row 21
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 4. Update period values (E16:E23) and Valor Mora values (F16:F23) ---
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56000

$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 56000

$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56000

$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 56000

$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 56000

$ws.Range("E21").Value = "2502"
$ws.Range("F21").Value = 56000

$ws.Range("E22").Value = "2501"
$ws.Range("F22").Value = 56000

$ws.Range("E23").Value = "2412"
$ws.Range("F23").Value = 24267

# Note: inserting the new data row at 23 above already shifted the old row 27
# ("___" separator) down to row 28, and the old row 28 ("NOMBRE.../FIRMA...")
# down to row 29 -- which is exactly the target layout. No further row
# insertion is required.

Write-Host "Done"
